# Add "Quelimane" as a new district row under the ZAMBEZIA province,
# right after the existing "Pebane" row (mentoring-core mapping-districts).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row goes right below the current last data row (row 16).
$newRow = 17

$ws.Range("B$newRow").Value = "ZAMBEZIA"
$ws.Range("C$newRow").Value = "Quelimane"

# Match the formatting used by the rest of the district (column C) entries,
# e.g. the cell directly above the new one.
$ws.Range("C$newRow").Style = $ws.Range("C16").Style

# The active cell ends up on the district column of the row above the new
# one, matching where the user would land after typing the new entry.
$ws.Range("C16").Select() | Out-Null
